$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2023-10-01 Sunday" "2023-10-02 Monday"

Replace-Text "31×70=2170" "48×35=1680"
Replace-Text "53×89=4717" "40×16=640"
Replace-Text "26×60=1560" "90×15=1350"
Replace-Text "24×89=2136" "81×24=1944"
Replace-Text "33×95=3135" "21×61=1281"

Replace-Text "98×12=1176" "48×50=2400"
Replace-Text "61×30=1830" "60×83=4980"
Replace-Text "75×87=6525" "39×22=858"
Replace-Text "18×91=1638" "87×67=5829"
Replace-Text "74×38=2812" "15×82=1230"

Replace-Text "53×93=4929" "90×71=6390"
Replace-Text "35×99=3465" "43×87=3741"
Replace-Text "72×64=4608" "15×11=165"
Replace-Text "48×11=528" "12×51=612"
Replace-Text "34×11=374" "18×24=432"

Replace-Text "59×88=5192" "63×89=5607"
Replace-Text "19×40=760" "91×34=3094"
Replace-Text "19×87=1653" "11×93=1023"
Replace-Text "29×57=1653" "43×24=1032"
Replace-Text "59×70=4130" "30×49=1470"

Replace-Text "28×97=2716" "64×59=3776"
Replace-Text "97×23=2231" "21×48=1008"
Replace-Text "72×37=2664" "97×27=2619"
Replace-Text "87×23=2001" "60×94=5640"
Replace-Text "26×14=364" "61×11=671"
